$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.383.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.16%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.046.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.15%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'228.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.36%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'55.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.40%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -2.77%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0807"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.93%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D12").Value = "'2.349.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'14.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -5.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -3.70%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.751"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.77%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -2.34%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.044.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.04%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'37.251.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.49%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.46%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'69.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.94%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0845"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.72%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'225.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.94%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.18%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E25").Value = "'  -5.20%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -3.21%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'168.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -4.24%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -1.17%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'18.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.38%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -2.89%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.07%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -3.80%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.53%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.11%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.07%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -5.04%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.44%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.502.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.45%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -7.14%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.94%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'16.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.28%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'95.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -5.53%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0930"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.52%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -4.05%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -4.96%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.72%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.96%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'FTXToken"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'3.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -9.65%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'RocketPoolETH"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'2.236.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.08%  "
$ws.Range("E51").Style = "Normal"
